# Removed cards from the Duel Decks
# The workbook contains an album/price list. Two rows referencing the
# "Duel Decks: Izzet vs Golgari" set ("Life from the Loam" and
# "Jarad, Golgari Lich Lord") are removed, the remaining rows shift up,
# and the SUM formula below the table is adjusted to the new range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 32 = "Life from the Loam" / "Duel Decks: Izzet vs Golgari"
$ws.Rows("32").Delete()

# After the first deletion, row 52 ("Jarad, Golgari Lich Lord" /
# "Duel Decks: Izzet vs Golgari") has shifted up to row 51.
$ws.Rows("51").Delete()

# Reset the view: scroll back to the top and select A51, matching the
# saved workbook state.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("A51").Select()
